# Update "想去人数" (want-to-go count) figures for the 展览 and 全部类型 sheets
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 2186
    $ws.Range("F5").Value = 188
}
